# Insert a new data row at row 55 (pushes existing rows 55..166 down to 56..167)
# and populate it with the new weekly record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 55, shifting data down and extending
# the used range from A1:R166 to A1:R167.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record.
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44775
$ws.Cells.Item(55, 4).Style = $ws.Cells.Item(56, 4).Style
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = 100112012
$ws.Cells.Item(55, 7).Value = "Espinaca"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 65
$ws.Cells.Item(55, 11).Value = 9000
$ws.Cells.Item(55, 12).Value = 9000
$ws.Cells.Item(55, 13).Value = 9000
$ws.Cells.Item(55, 14).Value = "$/docena de atados"
$ws.Cells.Item(55, 15).Value = "Región Metropolitana"
$ws.Cells.Item(55, 16).Value = 3000
$ws.Cells.Item(55, 17).Value = 3
$ws.Cells.Item(55, 18).Value = "Hortaliza"
